$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the election results for row 2 (SANTARÉM / CONSTÂNCIA)
$ws.Range("H2").Value = 47
$ws.Range("I2").Value = 92
$ws.Range("J2").Value = 401
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 135
$ws.Range("M2").Value = 6
$ws.Range("N2").Value = 63
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 7
$ws.Range("S2").Value = 44
$ws.Range("T2").Value = 73
$ws.Range("U2").Value = 4
$ws.Range("V2").Value = 669
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 628
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 7
